# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the profit sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2972.5264
$ws.Range("J86").Value = 4101.375
$ws.Range("L86").Value = 4101.375
$ws.Range("N86").Value = -6347.375
$ws.Range("H89").Value = 2972.5264
$ws.Range("J89").Value = 4101.375
$ws.Range("L89").Value = 20506.875
$ws.Range("N89").Value = -31738.875
$ws.Range("H100").Value = 2152.8333
$ws.Range("I100").Value = 2201.8
$ws.Range("J100").Value = 1908
$ws.Range("K100").Value = 2201.8
$ws.Range("L100").Value = 1908
$ws.Range("M100").Value = -1660.8
$ws.Range("N100").Value = -2990
$ws.Range("H116").Value = 7590.5713
$ws.Range("I116").Value = 5918.6665
$ws.Range("K116").Value = 5918.6665
$ws.Range("M116").Value = -2476.6665
$ws.Range("H131").Value = 13342.667
$ws.Range("J131").Value = 34182
$ws.Range("L131").Value = 102546
$ws.Range("N131").Value = -112626
$ws.Range("H137").Value = 1507.9
$ws.Range("I137").Value = 1451.7317
$ws.Range("K137").Value = 4355.1951
$ws.Range("M137").Value = -1805.1951
$ws.Range("H138").Value = 8067494.5
$ws.Range("I138").Value = 1485
$ws.Range("J138").Value = 11367226
$ws.Range("K138").Value = 4455
$ws.Range("L138").Value = 34101678
$ws.Range("M138").Value = 685
$ws.Range("N138").Value = -34111958
$ws.Range("H141").Value = 2346.0715
$ws.Range("I141").Value = 2346.0715
$ws.Range("K141").Value = 7038.2145
$ws.Range("M141").Value = -1858.2145

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5407.04
$ws.Range("I45").Value = 6887.778
$ws.Range("J45").Value = 1599.4286
$ws.Range("K45").Value = 6887.778
$ws.Range("L45").Value = 1599.4286
$ws.Range("M45").Value = -6510.778
$ws.Range("N45").Value = -2353.4286
$ws.Range("H74").Value = 6913.4863
$ws.Range("I74").Value = 1381.4286
$ws.Range("K74").Value = 1381.4286
$ws.Range("M74").Value = -507.4286
$ws.Range("H76").Value = 100000
$ws.Range("J76").Value = 100000
$ws.Range("L76").Value = 100000
$ws.Range("N76").Value = -100676
$ws.Range("H77").Value = 6913.4863
$ws.Range("I77").Value = 1381.4286
$ws.Range("K77").Value = 6907.143
$ws.Range("M77").Value = -2539.143
$ws.Range("H79").Value = 100000
$ws.Range("J79").Value = 100000
$ws.Range("L79").Value = 100000
$ws.Range("N79").Value = -102340
$ws.Range("H97").Value = 2170.3
$ws.Range("I97").Value = 1550.3334
$ws.Range("J97").Value = 7750
$ws.Range("K97").Value = 1550.3334
$ws.Range("L97").Value = 7750
$ws.Range("M97").Value = -1054.3334
$ws.Range("N97").Value = -8742
$ws.Range("H102").Value = 2025.375
$ws.Range("J102").Value = 1016
$ws.Range("L102").Value = 1016
$ws.Range("N102").Value = -4260
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 3360.5305
$ws.Range("I132").Value = 3016.0977
$ws.Range("K132").Value = 9048.293099999999
$ws.Range("M132").Value = -6518.293099999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4761.364
$ws.Range("I86").Value = 6458.231
$ws.Range("J86").Value = 2310.3333
$ws.Range("K86").Value = 6458.231
$ws.Range("L86").Value = 2310.3333
$ws.Range("M86").Value = -5335.231
$ws.Range("N86").Value = -4556.3333
$ws.Range("H89").Value = 4761.364
$ws.Range("I89").Value = 6458.231
$ws.Range("J89").Value = 2310.3333
$ws.Range("K89").Value = 32291.155
$ws.Range("L89").Value = 11551.6665
$ws.Range("M89").Value = -26675.155
$ws.Range("N89").Value = -22783.6665
$ws.Range("H99").Value = 148400.58
$ws.Range("I99").Value = 102760.8
$ws.Range("K99").Value = 102760.8
$ws.Range("M99").Value = -101262.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38228.168
$ws.Range("I31").Value = 47993.824
$ws.Range("J31").Value = 6141
$ws.Range("K31").Value = 47993.824
$ws.Range("L31").Value = 6141
$ws.Range("M31").Value = -47698.824
$ws.Range("N31").Value = -6731
$ws.Range("H34").Value = 38228.168
$ws.Range("I34").Value = 47993.824
$ws.Range("J34").Value = 6141
$ws.Range("K34").Value = 47993.824
$ws.Range("L34").Value = 6141
$ws.Range("M34").Value = -47791.824
$ws.Range("N34").Value = -6545
$ws.Range("H58").Value = 3038.1052
$ws.Range("I58").Value = 2151.5557
$ws.Range("J58").Value = 3836
$ws.Range("K58").Value = 2151.5557
$ws.Range("L58").Value = 3836
$ws.Range("M58").Value = -1948.5557
$ws.Range("N58").Value = -4242
$ws.Range("H136").Value = 3038.1052
$ws.Range("I136").Value = 2151.5557
$ws.Range("J136").Value = 3836
$ws.Range("K136").Value = 6454.6671
$ws.Range("L136").Value = 11508
$ws.Range("M136").Value = -3904.6671
$ws.Range("N136").Value = -16608
$ws.Range("H141").Value = 413984.16
$ws.Range("J141").Value = 413984.16
$ws.Range("L141").Value = 413984.16
$ws.Range("N141").Value = -424344.16

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 270.6
$ws.Range("I10").Value = 213.25
$ws.Range("K10").Value = 639.75
$ws.Range("M10").Value = -500.75
$ws.Range("H16").Value = 1346.4
$ws.Range("I16").Value = 950
$ws.Range("K16").Value = 2850
$ws.Range("M16").Value = -2677
$ws.Range("H113").Value = 4546289
$ws.Range("I113").Value = 5556331
$ws.Range("K113").Value = 16668993
$ws.Range("M113").Value = -16666823
$ws.Range("H122").Value = 1614.4584
$ws.Range("I122").Value = 1635.5714
$ws.Range("J122").Value = 1584.9
$ws.Range("K122").Value = 14720.1426
$ws.Range("L122").Value = 14264.1
$ws.Range("M122").Value = -12270.1426
$ws.Range("N122").Value = -19164.1
$ws.Range("H129").Value = 16500963
$ws.Range("I129").Value = 19800556
$ws.Range("K129").Value = 59401668
$ws.Range("M129").Value = -59396668
$ws.Range("H137").Value = 2638
$ws.Range("J137").Value = 3407.1428
$ws.Range("L137").Value = 10221.4284
$ws.Range("N137").Value = -20421.4284

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 9707.916999999999
$ws.Range("I57").Value = 7000
$ws.Range("J57").Value = 11642.143
$ws.Range("K57").Value = 7000
$ws.Range("L57").Value = 11642.143
$ws.Range("M57").Value = -6180
$ws.Range("N57").Value = -13282.143
$ws.Range("H96").Value = 46753.332
$ws.Range("J96").Value = 46753.332
$ws.Range("L96").Value = 46753.332
$ws.Range("N96").Value = -52245.332
$ws.Range("H102").Value = 62500756
$ws.Range("I102").Value = 706.38464
$ws.Range("K102").Value = 706.38464
$ws.Range("M102").Value = 915.61536

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3222.6155
$ws.Range("I22").Value = 3899
$ws.Range("J22").Value = 3166.25
$ws.Range("K22").Value = 3899
$ws.Range("L22").Value = 3166.25
$ws.Range("M22").Value = -3604
$ws.Range("N22").Value = -3756.25
$ws.Range("H27").Value = 3222.6155
$ws.Range("I27").Value = 3899
$ws.Range("J27").Value = 3166.25
$ws.Range("K27").Value = 3899
$ws.Range("L27").Value = 3166.25
$ws.Range("M27").Value = -3792
$ws.Range("N27").Value = -3380.25
$ws.Range("H132").Value = 3268.4285
$ws.Range("I132").Value = 3040.2188
$ws.Range("J132").Value = 3998.7
$ws.Range("K132").Value = 9120.6564
$ws.Range("L132").Value = 11996.1
$ws.Range("M132").Value = -6590.6564
$ws.Range("N132").Value = -17056.1

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 225004500
$ws.Range("I20").Value = 450000000
$ws.Range("J20").Value = 9000
$ws.Range("K20").Value = 450000000
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = -449999760
$ws.Range("N20").Value = -9480
$ws.Range("H32").Value = 14679
$ws.Range("I32").Value = 6363
$ws.Range("K32").Value = 6363
$ws.Range("M32").Value = -6046
$ws.Range("H107").Value = 25088.781
$ws.Range("I107").Value = 708.6177
$ws.Range("K107").Value = 2125.8531
$ws.Range("M107").Value = -205.8531000000003
$ws.Range("H122").Value = 1434.4546
$ws.Range("I122").Value = 1405.7805
$ws.Range("J122").Value = 1826.3334
$ws.Range("K122").Value = 4217.3415
$ws.Range("L122").Value = 5479.0002
$ws.Range("M122").Value = -1767.3415
$ws.Range("N122").Value = -10379.0002
$ws.Range("H132").Value = 2499.6758
$ws.Range("I132").Value = 2266.1562
$ws.Range("J132").Value = 3994.2
$ws.Range("K132").Value = 6798.4686
$ws.Range("L132").Value = 11982.6
$ws.Range("M132").Value = -4268.4686
$ws.Range("N132").Value = -17042.6
$ws.Range("H136").Value = 2776.3704
$ws.Range("I136").Value = 1662.0625
$ws.Range("J136").Value = 4397.1816
$ws.Range("K136").Value = 4986.1875
$ws.Range("L136").Value = 13191.5448
$ws.Range("M136").Value = -2436.1875
$ws.Range("N136").Value = -18291.5448

